$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = '안녕하십니까.'
$ws.Range("B2").Value = 'Greetings (Formal).'
$ws.Range("A3").Value = '수고하셨습니다.'
$ws.Range("B3").Value = 'Good job today. / Thank you for your effort.'
$ws.Range("A4").Value = '회의 시작하겠습니다.'
$ws.Range("B4").Value = 'Let''s start the meeting.'
$ws.Range("A5").Value = '의견 있으신가요?'
$ws.Range("B5").Value = 'Do you have any opinions?'
$ws.Range("A6").Value = '결재 부탁드립니다.'
$ws.Range("B6").Value = 'Please approve this.'
$ws.Range("A7").Value = '보고서 제출했습니다.'
$ws.Range("B7").Value = 'I submitted the report.'
$ws.Range("A8").Value = '이메일 확인 부탁드립니다.'
$ws.Range("B8").Value = 'Please check your email.'
$ws.Range("A9").Value = '오늘 회식 있나요?'
$ws.Range("B9").Value = 'Do we have a team dinner today?'
$ws.Range("A10").Value = '야근해야 할 것 같습니다.'
$ws.Range("B10").Value = 'I think I have to work overtime.'
$ws.Range("A11").Value = '출장 다녀오겠습니다.'
$ws.Range("B11").Value = 'I will go on a business trip.'
$ws.Range("A12").Value = '죄송하지만, 잠시 통화 가능하신가요?'
$ws.Range("B12").Value = 'Sorry, can you talk for a moment?'
$ws.Range("A14").Value = '일정 확인해보겠습니다.'
$ws.Range("B14").Value = 'I will check the schedule.'
$ws.Range("A15").Value = '협조해주셔서 감사합니다.'
$ws.Range("B15").Value = 'Thank you for your cooperation.'
$ws.Range("A16").Value = '마감 기한은 언제인가요?'
$ws.Range("B16").Value = 'When is the deadline?'
$ws.Range("A17").Value = '문제 없습니다.'
$ws.Range("B17").Value = 'No problem.'
$ws.Range("A18").Value = '검토 후 연락드리겠습니다.'
$ws.Range("B18").Value = 'I will contact you after review.'
$ws.Range("A19").Value = '명함 한 장 주시겠습니까?'
$ws.Range("B19").Value = 'Could you give me your business card?'
$ws.Range("A20").Value = '승진 축하드립니다.'
$ws.Range("B20").Value = 'Congratulations on your promotion.'
$ws.Range("A21").Value = '휴가 잘 다녀오세요.'
$ws.Range("B21").Value = 'Have a nice vacation.'
$ws.Range("A22").Value = '퇴근하겠습니다.'
$ws.Range("B22").Value = 'I''m leaving work now.'
$ws.Range("A23").Value = '내일 뵙겠습니다.'
$ws.Range("B23").Value = 'See you tomorrow.'
$ws.Range("A24").Value = '양해 부탁드립니다.'
$ws.Range("B24").Value = 'I ask for your understanding.'
$ws.Range("A25").Value = '최선을 다하겠습니다.'
$ws.Range("B25").Value = 'I will do my best.'
$ws.Range("A26").Value = '좋은 결과 기대하겠습니다.'
$ws.Range("B26").Value = 'I look forward to good results.'
$ws.Range("A27").Value = '거래처 미팅이 있습니다.'
$ws.Range("B27").Value = 'I have a meeting with a client.'
$ws.Range("A28").Value = '견적서 보내주세요.'
$ws.Range("B28").Value = 'Please send me the quotation.'
$ws.Range("A29").Value = '계약이 성사되었습니다.'
$ws.Range("B29").Value = 'The contract has been signed.'
$ws.Range("A30").Value = '프로젝트 진행 상황 보고해주세요.'
$ws.Range("B30").Value = 'Please report the project progress.'
$ws.Range("A31").Value = '수고하세요.'
$ws.Range("B31").Value = 'Keep up the good work. (Formal greeting)'
$ws.Range("A32").Value = '자료 준비되었습니다.'
$ws.Range("B32").Value = 'The materials are ready.'
$ws.Range("A33").Value = '언제 시간 되시나요?'
$ws.Range("B33").Value = 'When are you available?'
$ws.Range("A34").Value = '점심 식사 하셨습니까?'
$ws.Range("B34").Value = 'Have you had lunch?'
$ws.Range("A35").Value = '커피 한 잔 하시겠습니까?'
$ws.Range("B35").Value = 'Would you like a cup of coffee?'
$ws.Range("A36").Value = '잠시 자리 비우셨습니다.'
$ws.Range("B36").Value = 'He/She is away from the desk for a moment.'
$ws.Range("A37").Value = '메모 남겨드릴까요?'
$ws.Range("B37").Value = 'Shall I leave a message?'
$ws.Range("A38").Value = '팩스로 보내주세요.'
$ws.Range("B38").Value = 'Please send it by fax.'
$ws.Range("A39").Value = '프레젠테이션 준비 완료했습니다.'
$ws.Range("B39").Value = 'Presentation preparation is complete.'
$ws.Range("A40").Value = '예산안 검토 부탁드립니다.'
$ws.Range("B40").Value = 'Please review the budget proposal.'
$ws.Range("A41").Value = '다음 주 일정 잡겠습니다.'
$ws.Range("B41").Value = 'I will schedule it for next week.'
$ws.Range("A42").Value = '출근했습니다.'
$ws.Range("B42").Value = 'I have arrived at work.'
$ws.Range("A43").Value = '병가 내겠습니다.'
$ws.Range("B43").Value = 'I will take a sick leave.'
$ws.Range("A44").Value = '연차 쓰겠습니다.'
$ws.Range("B44").Value = 'I will use my annual leave.'
$ws.Range("A45").Value = '급한 일입니까?'
$ws.Range("B45").Value = 'Is it urgent?'
$ws.Range("A46").Value = '우선순위가 어떻게 되나요?'
$ws.Range("B46").Value = 'What is the priority?'
$ws.Range("A47").Value = '피드백 부탁드립니다.'
$ws.Range("B47").Value = 'Please give me feedback.'
$ws.Range("A48").Value = '수정해서 다시 보내겠습니다.'
$ws.Range("B48").Value = 'I will revise and resend it.'
$ws.Range("A49").Value = '성공적인 프로젝트였습니다.'
$ws.Range("B49").Value = 'It was a successful project.'
$ws.Range("A50").Value = '함께 일해서 즐거웠습니다.'
$ws.Range("B50").Value = 'It was a pleasure working with you.'
$ws.Range("A51").Value = '앞으로도 잘 부탁드립니다.'
$ws.Range("B51").Value = 'I look forward to working with you.'
$ws.Range("A52").Value = '안녕하십니까.'
$ws.Range("B52").Value = 'Hello (formal).'
$ws.Range("A53").Value = '처음 뵙겠습니다.'
$ws.Range("B53").Value = 'Nice to meet you (first time).'
$ws.Range("A54").Value = '잘 부탁드립니다.'
$ws.Range("B54").Value = 'I look forward to your cooperation.'
$ws.Range("A55").Value = '오랜만입니다.'
$ws.Range("B55").Value = 'Long time no see (formal).'
$ws.Range("A56").Value = '그동안 잘 지내셨습니까?'
$ws.Range("B56").Value = 'How have you been (formal)?'
$ws.Range("A57").Value = '소개해 드리겠습니다.'
$ws.Range("B57").Value = 'Let me introduce you.'
$ws.Range("A58").Value = '이쪽은 제 동료입니다.'
$ws.Range("B58").Value = 'This is my colleague.'
$ws.Range("A59").Value = '명함 교환하시죠.'
$ws.Range("B59").Value = 'Let''s exchange business cards.'
$ws.Range("A60").Value = '연락처를 알 수 있을까요?'
$ws.Range("B60").Value = 'May I have your contact information?'
$ws.Range("A61").Value = '이메일 주소 알려주세요.'
$ws.Range("B61").Value = 'Please tell me the email address.'
$ws.Range("A62").Value = '전화 연결해 드리겠습니다.'
$ws.Range("B62").Value = 'I''ll connect you.'
$ws.Range("A63").Value = '부재중입니다.'
$ws.Range("B63").Value = 'He/She is currently unavailable/out of office.'
$ws.Range("A64").Value = '나중에 다시 걸겠습니다.'
$ws.Range("B64").Value = 'I''ll call back later.'
$ws.Range("A65").Value = '메시지를 남기시겠습니까?'
$ws.Range("B65").Value = 'Would you like to leave a message?'
$ws.Range("A66").Value = '잘 안 들립니다.'
$ws.Range("B66").Value = 'I can''t hear you well.'
$ws.Range("A67").Value = '좀 더 크게 말씀해 주세요.'
$ws.Range("B67").Value = 'Please speak up a little.'
$ws.Range("A68").Value = '끊겠습니다.'
$ws.Range("B68").Value = 'I''m hanging up.'
$ws.Range("A69").Value = '회의실을 예약했습니다.'
$ws.Range("B69").Value = 'I reserved the meeting room.'
$ws.Range("A70").Value = '빔 프로젝터가 필요합니다.'
$ws.Range("B70").Value = 'I need a beam projector.'
$ws.Range("A71").Value = '자료를 복사해 주세요.'
$ws.Range("B71").Value = 'Please copy the materials.'
$ws.Range("A72").Value = '회의록을 작성하겠습니다.'
$ws.Range("B72").Value = 'I will write the minutes.'
$ws.Range("A73").Value = '안건이 무엇입니까?'
$ws.Range("B73").Value = 'What is the agenda?'
$ws.Range("A74").Value = '결론을 내립시다.'
$ws.Range("B74").Value = 'Let''s come to a conclusion.'
$ws.Range("A75").Value = '다음 회의는 언제입니까?'
$ws.Range("B75").Value = 'When is the next meeting?'
$ws.Range("A76").Value = '참석해 주셔서 감사합니다.'
$ws.Range("B76").Value = 'Thank you for attending.'
$ws.Range("A77").Value = '추가 질문 있습니까?'
$ws.Range("B77").Value = 'Any further questions?'
$ws.Range("A78").Value = '시간이 부족합니다.'
$ws.Range("B78").Value = 'We are running out of time.'
$ws.Range("A79").Value = '잠시 쉬겠습니다.'
$ws.Range("B79").Value = 'Let''s take a break.'
$ws.Range("A80").Value = '이만 마치겠습니다.'
$ws.Range("B80").Value = 'We will wrap up now.'
$ws.Range("A81").Value = '보고서를 검토했습니다.'
$ws.Range("B81").Value = 'I reviewed the report.'
$ws.Range("A82").Value = '수정이 필요합니다.'
$ws.Range("B82").Value = 'It needs revision.'
$ws.Range("A83").Value = '오타가 있습니다.'
$ws.Range("B83").Value = 'There is a typo.'
$ws.Range("A84").Value = '내용이 부족합니다.'
$ws.Range("B84").Value = 'The content is insufficient.'
$ws.Range("A85").Value = '다시 작성해 주세요.'
$ws.Range("B85").Value = 'Please rewrite it.'
$ws.Range("A86").Value = '잘 작성되었습니다.'
$ws.Range("B86").Value = 'It is well written.'
$ws.Range("A87").Value = '승인 받았습니다.'
$ws.Range("B87").Value = 'It has been approved.'
$ws.Range("A88").Value = '반려되었습니다.'
$ws.Range("B88").Value = 'It has been rejected.'
$ws.Range("A89").Value = '기한을 연장할 수 있나요?'
$ws.Range("B89").Value = 'Can we extend the deadline?'
$ws.Range("A90").Value = '최대한 빨리 처리하겠습니다.'
$ws.Range("B90").Value = 'I will process it as soon as possible.'
$ws.Range("A91").Value = '문제가 해결되었습니다.'
$ws.Range("B91").Value = 'The problem has been resolved.'
$ws.Range("A92").Value = '컴퓨터가 고장났습니다.'
$ws.Range("B92").Value = 'The computer is broken.'
$ws.Range("A93").Value = '인터넷이 안 됩니다.'
$ws.Range("B93").Value = 'The internet is not working.'
$ws.Range("A94").Value = '비밀번호를 잊어버렸습니다.'
$ws.Range("B94").Value = 'I forgot my password.'
$ws.Range("A95").Value = '로그인이 안 됩니다.'
$ws.Range("B95").Value = 'I can''t log in.'
$ws.Range("A96").Value = '파일이 안 열립니다.'
$ws.Range("B96").Value = 'The file won''t open.'
$ws.Range("A97").Value = '바이러스에 감염되었습니다.'
$ws.Range("B97").Value = 'It is infected with a virus.'
$ws.Range("A98").Value = '백업을 해야 합니다.'
$ws.Range("B98").Value = 'We need to backup.'
$ws.Range("A99").Value = '소프트웨어를 업데이트하세요.'
$ws.Range("B99").Value = 'Please update the software.'
$ws.Range("A100").Value = '재부팅해 보세요.'
$ws.Range("B100").Value = 'Try rebooting.'
$ws.Range("A101").Value = 'IT 지원팀에 연락하세요.'
$ws.Range("B101").Value = 'Contact IT support.'
